# Applies the "Doc de requerimientos completado y sumamos estructura front y back"
# edit:
#   1. Removes the stray "_GoBack" bookmark that sat at the very start of the
#      document (right after the title paragraph's pPr).
#   2. Fills in the trailing blank paragraph with "OPCIONAL:".
#   3. Turns the final (completely empty) paragraph into a real paragraph with
#      the "Historial de ventas..." sentence, and re-homes the "_GoBack"
#      bookmark there (now wrapping the end of that new text), matching where
#      Word leaves it after the last edit made to the document.

$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Drop the old _GoBack bookmark near the top of the document. -------
$d.Bookmarks.Item("_GoBack").Delete()

# --- 2. "OPCIONAL:" paragraph (currently blank) -----------------------------
$count = $d.Paragraphs.Count
$pOpcional = $d.Paragraphs.Item($count - 1)
$opcionalXml = '<w:p ' + $w + '>' +
    '<w:pPr>' +
        '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w:lang w:eastAsia="es-AR"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w:lang w:eastAsia="es-AR"/>' +
        '</w:rPr>' +
        '<w:t>OPCIONAL:</w:t>' +
    '</w:r>' +
'</w:p>'
[void]$pOpcional.Range.InsertXML($opcionalXml)

# --- 3. Final paragraph: new sentence + the relocated _GoBack bookmark. ----
$pFinal = $d.Paragraphs.Item($count)
$finalXml = '<w:p ' + $w + '>' +
    '<w:pPr>' +
        '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w:lang w:eastAsia="es-AR"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w:lang w:eastAsia="es-AR"/>' +
        '</w:rPr>' +
        '<w:t>Historial de ventas que muestre la &#8220;factura&#8221; con lo que haya comprado, accedible.</w:t>' +
    '</w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
'</w:p>'
[void]$pFinal.Range.InsertXML($finalXml)
